# dithionite_sensitivity/parameters.xlsx -- "Updates to mads sensitivity"
#
# The real, data-level change in this commit is the base value of the
# k_s2o4_disp parameter on the "summary" sheet (cell B2): it is tightened
# from 3.61E-005 to 1E-005. Everything else in the workbook (the mads and
# mads_tightened sheets) references summary!B2 through formulas, so those
# dependent cells recompute automatically once B2 changes.
#
# The commit also records updated cursor/selection positions left behind
# on the "summary" and "mads_tightened" sheets.

$wb = $excel.ActiveWorkbook

$wsSummary        = $wb.Worksheets.Item("summary")
$wsMadsTightened  = $wb.Worksheets.Item("mads_tightened")

# --- Core data edit -------------------------------------------------------
# summary!B2 (base value for k_s2o4_disp): 3.61E-005 -> 1E-005
$wsSummary.Range("B2").Value = 0.00001

# --- Leftover selection/cursor state recorded in the saved file ----------
# summary sheet: selection moves from B6 to B3
$wsSummary.Activate()
$wsSummary.Range("B3").Select()

# mads_tightened sheet: selection moves from C16 to E24
$wsMadsTightened.Activate()
$wsMadsTightened.Range("E24").Select()
